$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: город Котлас, female, 2019
$ws.Range("A32").Value = 11710000
$ws.Range("B32").Value = "город Котлас"
$ws.Range("C32").Value = "female"
$ws.Range("D32").Value = 2019
$ws.Range("E32").Value = 0.0588
$ws.Range("F32").Value = 0.0678
$ws.Range("G32").Value = 0.0604
$ws.Range("H32").Value = 0.0469
$ws.Range("I32").Value = 0.03168
$ws.Range("J32").Value = 0.04794
$ws.Range("K32").Value = 0.07715
$ws.Range("L32").Value = 0.08655
$ws.Range("M32").Value = 0.08026
$ws.Range("N32").Value = 0.06305
$ws.Range("O32").Value = 0.1276
$ws.Range("P32").Value = 0.0804
$ws.Range("Q32").Value = 0.0937
$ws.Range("R32").Value = 0.0779

# Row 33: город Котлас, male, 2019
$ws.Range("A33").Value = 11710000
$ws.Range("B33").Value = "город Котлас"
$ws.Range("C33").Value = "male"
$ws.Range("D33").Value = 2019
$ws.Range("E33").Value = 0.07043
$ws.Range("F33").Value = 0.08215
$ws.Range("G33").Value = 0.06976
$ws.Range("H33").Value = 0.0541
$ws.Range("I33").Value = 0.02719
$ws.Range("J33").Value = 0.04855
$ws.Range("K33").Value = 0.088
$ws.Range("L33").Value = 0.1001
$ws.Range("M33").Value = 0.09454
$ws.Range("N33").Value = 0.0683
$ws.Range("O33").Value = 0.11115
$ws.Range("P33").Value = 0.0675
$ws.Range("Q33").Value = 0.0665
$ws.Range("R33").Value = 0.05167

# Apply the same style as the surrounding data rows (style index 1 == no special format)
$ws.Range("A31:R31").Copy()
$ws.Range("A32:R33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the view: scroll to A13, select B34 (matches the committed sheetView state)
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B34").Select()
